$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price values that look like plain numbers (single
# decimal point). Excel auto-converts such text to a Number type on assignment,
# which would corrupt values such as "18.10" (-> 18.1) or change the stored
# cell type. Force these specific cells to Text format first so the exact
# string is preserved, matching how the workbook already stores these columns.
$textFormatCells = @(
    "D5",
    "D6",
    "D10",
    "D16",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D29",
    "D30",
    "D33",
    "D34",
    "D37",
    "D38",
    "D39",
    "D40",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D51"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "58.154.84"
$ws.Cells.Item(2, 5).Value = "  +0.52%  "
$ws.Cells.Item(3, 4).Value = "2.463.82"
$ws.Cells.Item(3, 5).Value = "  +0.54%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Value = "510.92"
$ws.Cells.Item(5, 5).Value = "  -2.50%  "
$ws.Cells.Item(6, 4).Value = "134.07"
$ws.Cells.Item(6, 5).Value = "  +3.37%  "
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "2.463.47"
$ws.Cells.Item(9, 5).Value = "  +0.37%  "
$ws.Cells.Item(10, 4).Value = "0.0984"
$ws.Cells.Item(10, 5).Value = "  +0.81%  "
$ws.Cells.Item(11, 5).Value = "  -0.56%  "
$ws.Cells.Item(12, 5).Value = "  +1.04%  "
$ws.Cells.Item(13, 5).Value = "  -6.02%  "
$ws.Cells.Item(14, 4).Value = "2.899.67"
$ws.Cells.Item(14, 5).Value = "  +0.52%  "
$ws.Cells.Item(15, 4).Value = "58.020.52"
$ws.Cells.Item(15, 5).Value = "  +0.42%  "
$ws.Cells.Item(16, 4).Value = "22.02"
$ws.Cells.Item(16, 5).Value = "  +2.23%  "
$ws.Cells.Item(17, 5).Value = "  +1.94%  "
$ws.Cells.Item(18, 4).Value = "2.430.40"
$ws.Cells.Item(18, 5).Value = "  -0.95%  "
$ws.Cells.Item(19, 4).Value = "10.39"
$ws.Cells.Item(19, 5).Value = "  +0.00%  "
$ws.Cells.Item(20, 4).Value = "4.13"
$ws.Cells.Item(20, 5).Value = "  +0.28%  "
$ws.Cells.Item(21, 4).Value = "315.20"
$ws.Cells.Item(21, 5).Value = "  +0.80%  "
$ws.Cells.Item(22, 4).Value = "6.48"
$ws.Cells.Item(22, 5).Value = "  +6.16%  "
$ws.Cells.Item(23, 4).Value = "0.999"
$ws.Cells.Item(23, 5).Value = "  -0.05%  "
$ws.Cells.Item(24, 4).Value = "5.72"
$ws.Cells.Item(24, 5).Value = "  -1.83%  "
$ws.Cells.Item(25, 4).Value = "65.48"
$ws.Cells.Item(25, 5).Value = "  +0.87%  "
$ws.Cells.Item(26, 4).Value = "0.997"
$ws.Cells.Item(26, 5).Value = "  -0.35%  "
$ws.Cells.Item(27, 5).Value = "  -0.19%  "
$ws.Cells.Item(28, 5).Value = "  -4.87%  "
$ws.Cells.Item(29, 4).Value = "7.63"
$ws.Cells.Item(29, 5).Value = "  +5.33%  "
$ws.Cells.Item(30, 4).Value = "171.68"
$ws.Cells.Item(30, 5).Value = "  -1.24%  "
$ws.Cells.Item(31, 4).Value = "0.0₃0739"
$ws.Cells.Item(31, 5).Value = "  +0.68%  "
$ws.Cells.Item(32, 5).Value = "  +0.36%  "
$ws.Cells.Item(33, 4).Value = "6.17"
$ws.Cells.Item(33, 5).Value = "  +0.31%  "
$ws.Cells.Item(34, 4).Value = "1.15"
$ws.Cells.Item(34, 5).Value = "  +1.23%  "
$ws.Cells.Item(35, 5).Value = "  +0.07%  "
$ws.Cells.Item(36, 5).Value = "  +0.14%  "
$ws.Cells.Item(37, 4).Value = "18.10"
$ws.Cells.Item(37, 5).Value = "  +1.65%  "
$ws.Cells.Item(38, 4).Value = "1.25"
$ws.Cells.Item(38, 5).Value = "  +5.86%  "
$ws.Cells.Item(39, 4).Value = "3.92"
$ws.Cells.Item(39, 5).Value = "  +3.88%  "
$ws.Cells.Item(40, 4).Value = "36.83"
$ws.Cells.Item(40, 5).Value = "  +1.53%  "
$ws.Cells.Item(41, 5).Value = "  +2.29%  "
$ws.Cells.Item(42, 5).Value = "  +0.81%  "
$ws.Cells.Item(43, 4).Value = "136.99"
$ws.Cells.Item(43, 5).Value = "  +11.57%  "
$ws.Cells.Item(44, 4).Value = "3.42"
$ws.Cells.Item(44, 5).Value = "  +0.91%  "
$ws.Cells.Item(45, 4).Value = "4.93"
$ws.Cells.Item(45, 5).Value = "  +3.08%  "
$ws.Cells.Item(46, 5).Value = "  -0.90%  "
$ws.Cells.Item(47, 4).Value = "256.04"
$ws.Cells.Item(47, 5).Value = "  -0.24%  "
$ws.Cells.Item(48, 4).Value = "0.0919"
$ws.Cells.Item(49, 5).Value = "  +0.86%  "
$ws.Cells.Item(50, 5).Value = "  +2.38%  "
$ws.Cells.Item(51, 4).Value = "17.31"
$ws.Cells.Item(51, 5).Value = "  +1.95%  "
